$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("F").Insert()
$ws.Range("F1").Value = "ServiceTeam"
$ws.Range("F2").Value = "a1Nq0000000RlgV"
